$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 "description": update note text on B1 -------------------------
$ws2.Range("B1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- Sheet2: add two helper-note columns (G, H) at the top ----------------
# Row1: blank, centered cells above the merged note row
$ws2.Range("G1:H1").Style = "Normal"
$ws2.Range("G1:H1").HorizontalAlignment = -4108

# Row2: labelled notes ("can be blank" / "must not be blank") using the
# same Neutral / Bad cell styles already used elsewhere on the sheet
$ws2.Range("G2").Value = "เป็นค่าว่างได้"
$ws2.Range("B2").Copy()
$ws2.Range("G2").PasteSpecial(-4122)

$ws2.Range("H2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"
$ws2.Range("A2").Copy()
$ws2.Range("H2").PasteSpecial(-4122)

# Merge G1:H1 into a single note cell
$ws2.Range("G1:H1").Merge()

# Column widths for the two new columns
$ws2.Range("G1:H1").EntireColumn.ColumnWidth = 25.7265625

# --- Sheet2 view/selection -------------------------------------------------
$ws2.Activate()
$ws2.Range("G1:H1048576").Select()
$excel.ActiveWindow.ScrollColumn = 2

# --- Sheet1: selection moves back to A2 ------------------------------------
$ws1.Activate()
$ws1.Range("A2").Select()
